$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-19 Tuesday", "2025-08-20 Wednesday"),
    @("478÷3=", "480÷4="),
    @("720÷8=", "655÷4="),
    @("844÷3=", "123÷6="),
    @("928÷9=", "496÷4="),
    @("396÷3=", "348÷9="),
    @("735÷9=", "850÷5="),
    @("599÷5=", "867÷5="),
    @("825÷7=", "722÷2="),
    @("761÷5=", "448÷6="),
    @("799÷8=", "963÷3="),
    @("316÷9=", "958÷7="),
    @("863÷9=", "330÷4="),
    @("234÷4=", "344÷6="),
    @("627÷9=", "776÷4="),
    @("979÷7=", "630÷9="),
    @("761÷4=", "911÷4="),
    @("419÷2=", "264÷9="),
    @("514÷5=", "832÷6="),
    @("593÷9=", "464÷9="),
    @("676÷2=", "723÷3="),
    @("395÷7=", "243÷4="),
    @("480÷5=", "630÷4="),
    @("802÷5=", "138÷2="),
    @("241÷6=", "670÷4="),
    @("615÷6=", "232÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
